# Prise en compte des devoirs dans la range de deadline : la liste des
# "devoirs" (Activités) pour Sloth/Skido a été rafraîchie, et la colonne
# dupliquée "Activité 112" (2e devoir du même nom qu'"Activité 11") est
# retirée de la range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Preserve the per-status cell formatting (fill colours) before the
#     values move around ---
# E2 ("Reçu", orange/EBBD86 fill) becomes the new B2.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("C2:D2").PasteSpecial(-4122) | Out-Null         # xlPasteFormats
$excel.CutCopyMode = 0

# B2's current fill (FFDFB3, used for timestamp cells) is reused for the
# two new timestamp cells C2/D2.
$ws.Range("E2").Copy() | Out-Null
$ws.Range("B2").PasteSpecial(-4122) | Out-Null            # xlPasteFormats
$excel.CutCopyMode = 0

# --- Update header row with the (now non-duplicated) homework names ---
$ws.Range("B1").Value = "Activité 23"
$ws.Range("C1").Value = "Activité 111"
$ws.Range("D1").Value = "Activité 11"

# --- Refresh Sloth's statuses for the remaining homeworks ---
$ws.Range("B2").Value = "Reçu"
$ws.Range("C2").Value = "le 23/05 à 01:46"
$ws.Range("D2").Value = "le 23/05 à 01:45"

# Skido's row stays "Non rendu" across the board.
$ws.Range("B3").Value = "Non rendu"
$ws.Range("C3").Value = "Non rendu"
$ws.Range("D3").Value = "Non rendu"

# --- Drop the obsolete duplicate "Activité 112" column (E) ---
$ws.Range("E1:E3").Clear() | Out-Null
